$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (F column) values for rows 6, 8, 13, 15 - repull data / recalc
$ws.Range("F6").Value = -10
$ws.Range("F8").Value = -4
$ws.Range("F13").Value = -10
$ws.Range("F15").Value = 2
